$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($ws, $cellRef, $val) {
    $rng = $ws.Range($cellRef)
    $rng.Value = "'" + $val
    $rng.Style = "Normal"
}

Set-TextValue $ws 'D2' '72.663.74'
Set-TextValue $ws 'E2' '  +0.97%  '
Set-TextValue $ws 'D3' '3.970.00'
Set-TextValue $ws 'E3' '  -1.28%  '
Set-TextValue $ws 'E4' '  +0.06%  '
Set-TextValue $ws 'D5' '618.28'
Set-TextValue $ws 'E5' '  +15.20%  '
Set-TextValue $ws 'D6' '166.18'
Set-TextValue $ws 'E6' '  +10.16%  '
Set-TextValue $ws 'D7' '0.685'
Set-TextValue $ws 'E7' '  -0.53%  '
Set-TextValue $ws 'E8' '  +0.03%  '
Set-TextValue $ws 'D9' '0.758'
Set-TextValue $ws 'E9' '  +0.66%  '
Set-TextValue $ws 'E10' '  -1.84%  '
Set-TextValue $ws 'D11' '57.86'
Set-TextValue $ws 'E11' '  +6.39%  '
Set-TextValue $ws 'D12' '0.0000314'
Set-TextValue $ws 'E12' '  -2.86%  '
Set-TextValue $ws 'D13' '11.26'
Set-TextValue $ws 'E13' '  +4.26%  '
Set-TextValue $ws 'D14' '4.603.58'
Set-TextValue $ws 'E14' '  -0.40%  '
Set-TextValue $ws 'D15' '3.971.83'
Set-TextValue $ws 'E15' '  -0.39%  '
Set-TextValue $ws 'D16' '1.27'
Set-TextValue $ws 'E16' '  +6.49%  '
Set-TextValue $ws 'D17' '14.21'
Set-TextValue $ws 'E17' '  +0.52%  '
Set-TextValue $ws 'D18' '20.71'
Set-TextValue $ws 'E18' '  +0.02%  '
Set-TextValue $ws 'E19' '  +0.28%  '
Set-TextValue $ws 'D20' '72.568.39'
Set-TextValue $ws 'E20' '  +1.04%  '
Set-TextValue $ws 'D21' '439.38'
Set-TextValue $ws 'E21' '  +1.21%  '
Set-TextValue $ws 'D22' '4.94'
Set-TextValue $ws 'E22' '  +17.22%  '
Set-TextValue $ws 'D23' '96.47'
Set-TextValue $ws 'E23' '  -1.70%  '
Set-TextValue $ws 'D24' '3.42'
Set-TextValue $ws 'E24' '  -3.79%  '
Set-TextValue $ws 'D25' '14.59'
Set-TextValue $ws 'E25' '  -0.70%  '
Set-TextValue $ws 'D26' '4.26'
Set-TextValue $ws 'E26' '  -0.76%  '
Set-TextValue $ws 'D27' '11.31'
Set-TextValue $ws 'E27' '  -2.22%  '
Set-TextValue $ws 'D28' '10.59'
Set-TextValue $ws 'E28' '  -2.59%  '
Set-TextValue $ws 'E29' '  +0.64%  '
Set-TextValue $ws 'D30' '36.16'
Set-TextValue $ws 'E30' '  -2.26%  '
Set-TextValue $ws 'D31' '7.83'
Set-TextValue $ws 'E31' '  -6.07%  '
Set-TextValue $ws 'D32' '13.96'
Set-TextValue $ws 'E32' '  +3.04%  '
Set-TextValue $ws 'D33' '0.132'
Set-TextValue $ws 'E33' '  -3.18%  '
Set-TextValue $ws 'B34' 'InjectiveProtocol'
Set-TextValue $ws 'C34' 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
Set-TextValue $ws 'D34' '48.38'
Set-TextValue $ws 'E34' '  -4.10%  '
Set-TextValue $ws 'B35' 'OKB'
Set-TextValue $ws 'C35' 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
Set-TextValue $ws 'D35' '72.08'
Set-TextValue $ws 'E35' '  +9.63%  '
Set-TextValue $ws 'D36' '641.46'
Set-TextValue $ws 'E36' '  -5.73%  '
Set-TextValue $ws 'E37' '  +7.55%  '
Set-TextValue $ws 'D38' '0.435'
Set-TextValue $ws 'E38' '  -5.11%  '
Set-TextValue $ws 'D39' '3.44'
Set-TextValue $ws 'E39' '  +1.90%  '
Set-TextValue $ws 'D40' '3.41'
Set-TextValue $ws 'E40' '  +1.36%  '
Set-TextValue $ws 'E41' '  -1.18%  '
Set-TextValue $ws 'E42' '  -0.20%  '
Set-TextValue $ws 'E43' '  +0.41%  '
Set-TextValue $ws 'D44' '0.0488'
Set-TextValue $ws 'E44' '  -0.98%  '
Set-TextValue $ws 'D45' '10.75'
Set-TextValue $ws 'E45' '  -1.16%  '
Set-TextValue $ws 'E46' '  +0.17%  '
Set-TextValue $ws 'B47' 'ApeXProtocol'
Set-TextValue $ws 'C47' 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
Set-TextValue $ws 'D47' '3.47'
Set-TextValue $ws 'E47' '  +3.63%  '
Set-TextValue $ws 'B48' 'Fetch.AI'
Set-TextValue $ws 'C48' 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
Set-TextValue $ws 'D48' '2.66'
Set-TextValue $ws 'E48' '  -0.64%  '
Set-TextValue $ws 'E49' '  -0.35%  '
Set-TextValue $ws 'D50' '2.903.50'
Set-TextValue $ws 'E50' '  +7.31%  '
Set-TextValue $ws 'B51' 'Monero'
Set-TextValue $ws 'C51' 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue $ws 'D51' '150.17'
Set-TextValue $ws 'E51' '  +4.09%  '
